# Update the workbook to add a new "Legal Project Manager" rate column (column J)
# to each of the four "Lot" worksheets.
#
# For every worksheet:
#   - Column J header (row 1) = "Legal Project Manager"
#   - Column J sub-header (row 2) = "Hourly rate" (same text as column I)
#   - Column J data rows = column I value * 3.25 (rounded the same way Excel would)
#   - Column J formatting/styles mirror column I (currency style, header styles, etc.)
#   - Column J width is set to match the other rate columns (just slightly narrower)

$wb = $excel.ActiveWorkbook

$sheetNames = @("Lot 1", "Lot 2", "Lot 3", "Lot 5")
$lastDataRow = @{ "Lot 1" = 23; "Lot 2" = 24; "Lot 3" = 24; "Lot 5" = 23 }

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Copy all formatting (styles) from column I onto column J first, covering
    # the header rows and every row that currently carries row-level formatting.
    $ws.Range("I1:I32").Copy()
    $ws.Range("J1").PasteSpecial(-4122)

    # Header text for the new column.
    $ws.Range("J1").Value = "Legal Project Manager"
    $ws.Range("J2").Value = "Hourly rate"

    # Fill in the new rate values (3.25x the Paralegal / Legal Assistant rate
    # found in column I) for every supplier row.
    $last = $lastDataRow[$name]
    for ($r = 3; $r -le $last; $r++) {
        $iValue = $ws.Cells.Item($r, 9).Value()
        $ws.Cells.Item($r, 10).Value = $iValue * 3.25
    }

    # Match the column width used elsewhere in the sheet for the rate columns.
    $ws.Columns.Item(10).ColumnWidth = 19.666666666666668
}

# Reset the "current selection" on the non-active lot sheets so that they no
# longer point at a stale cell/range (the sheets had leftover selections from
# when the rate tables ended at column I).
$ws1 = $wb.Worksheets.Item("Lot 1")
$ws2 = $wb.Worksheets.Item("Lot 2")
$ws3 = $wb.Worksheets.Item("Lot 3")
$ws4 = $wb.Worksheets.Item("Lot 5")

$ws2.Activate()
$ws2.Range("A1").Select()
$ws3.Activate()
$ws3.Range("A1").Select()
$ws4.Activate()
$ws4.Range("A1").Select()

$ws1.Activate()
$ws1.Range("A1").Select()
